$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Uniquer donors count" label -> "Unique donors count"
# (cell B7 on the "Report" sheet). Re-writing the cell text causes the
# shared-strings table to be rebuilt with the corrected string relocated
# to the end of the table, matching the authored change.
$ws.Range("B7").Value = "Unique donors count"

# Restore the default top-left scroll position and move the active
# selection to F6, matching the saved view state in the workbook.
$ws.Range("F6").Select()
